# Add a new row for LeetCode problem 57 "Insert Interval" to the tracking
# sheet, and refresh the view/window state to match the author's session.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data row (row 66), mirroring the layout of the existing rows:
# A=#, B=Name, C=Tags, D=Difficulty, E=Success, F=Fail, G=Time, H=First, I=Last Update
$row = 66
$ws.Cells.Item($row, 1).Value = 57
$ws.Cells.Item($row, 2).Value = "Insert Interval"
$ws.Cells.Item($row, 3).Value = "#array"
$ws.Cells.Item($row, 4).Value = "medium"
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 3
$ws.Cells.Item($row, 7).Value = 50
$ws.Cells.Item($row, 8).Value = 45847
$ws.Cells.Item($row, 9).Value = 45847

# Match styling of the date columns (H/I) used elsewhere in the sheet.
$ws.Range("H66").Style = $ws.Range("H65").Style
$ws.Range("I66").Style = $ws.Range("I65").Style
$ws.Range("A66").Style = $ws.Range("A65").Style
$ws.Range("D66").Style = $ws.Range("D65").Style
$ws.Range("E66").Style = $ws.Range("E65").Style
$ws.Range("F66").Style = $ws.Range("F65").Style
$ws.Range("G66").Style = $ws.Range("G65").Style
$ws.Range("B66").Style = $ws.Range("B65").Style
$ws.Range("C66").Style = $ws.Range("C65").Style

$ws.Rows.Item(66).RowHeight = 17

# Update the view: scroll back to column A and select the new "First"/"Last
# Update" cells for row 66.
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("H66:I66").Select()

# Widen the workbook window slightly, as recorded in the saved session.
$excel.Width = 29100
